$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates ---
# Add "(minimal install)" to the management-box TODO note.
$ws.Range("A17").Value = "NOTE: Management box - install Desktop Ubuntu GUI (minimal install) + go + vscode to debug deployer"

# --- Highlight the two edited "boxes needed" cells with a yellow fill ---
$ws.Range("D9:D10").Interior.Color = 65535

# --- Resource-spread numeric edits (prox1 table, rows 8-15) ---
$ws.Range("J8").Value = 174
$ws.Range("I9").Value = 40
$ws.Range("K9").ClearContents()
$ws.Range("I10").Value = 40
$ws.Range("I11").Value = 40
$ws.Range("I12").Value = 30
$ws.Range("K13").Value = 200
$ws.Range("K14").Value = 200
$ws.Range("K15").Value = 200

# --- Resource-spread numeric edits (prox2 table, row 23) ---
$ws.Range("J23").Value = 174

# --- Update selection to D10 ---
$ws.Range("D10").Select()

Write-Output "Applied deployer stack resource-spread updates"
